$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate the new Change Request row (row 13) with its values
$ws.Cells.Item(13, 1).Value = 12
$ws.Cells.Item(13, 2).Value = "Course details read from JSON file."
$ws.Cells.Item(13, 3).Value = "Dynamically load the Courses and course details against course.`nIn future if we want to add new course then we can add in JSON alone."
$ws.Cells.Item(13, 5).Value = "Completed"

# Copy the formatting (borders, alignment, wrap text, fonts) from the row above
$ws.Range("A12:E12").Copy()
$ws.Range("A13:E13").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Match the taller row height used for this multi-line entry
$ws.Rows.Item(13).RowHeight = 60

# Update the active cell/selection to the newly added cell
$ws.Range("B13").Select()
